# The "Maximum number of strata landed units" section used to show an
# inline picture (a screenshot of the intensity diagram). The picture is
# replaced by a plain hyperlink whose visible text is the image's URL.

$d = $word.ActiveDocument

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/SL02_Intensity.jpg?h=100%25&w=100%25"

# Find the inline picture that illustrates "Allowable intensity for a
# strata landed development" (there is exactly one in this document).
$target = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    if ($shape.AlternativeText -eq "Allowable intensity for a strata landed development") {
        $target = $shape
        break
    }
}

# Fall back to the first/only inline picture if the description ever
# drifts, so the edit still fires.
if (($target -eq $null) -and ($d.InlineShapes.Count -gt 0)) {
    $target = $d.InlineShapes.Item(1)
}

if ($target -ne $null) {
    $insertionRange = $target.Range
    $target.Delete()

    $link = $d.Hyperlinks.Add($insertionRange, $url, "", "", $url)
}
